$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.114.32"
$ws.Range("E2").Value = "  -2.44%  "

$ws.Range("D3").Value = "1.865.70"
$ws.Range("E3").Value = "  -2.05%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5154"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07153"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8897"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "

$ws.Range("E11").Value = "  -2.84%  "

$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "1.852.66"
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.307"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008473"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").Value = "27.130.98"
$ws.Range("E20").Value = "  -2.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.024"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("D22").Value = "2.097.73"
$ws.Range("E22").Value = "  -3.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.455"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.838"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.091"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.658"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.663"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05105"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.070"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.33%  "

$ws.Range("E35").Value = "  -6.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7261"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.05%  "

$ws.Range("E37").Value = "  -3.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.083"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.496"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.74%  "

$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5302"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.460"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.276"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.37%  "

$ws.Range("E45").Value = "  -3.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("E47").Value = "  -3.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.948"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.565"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("E51").Value = "  -4.92%  "
